{"js": "// Apply the two-digit-division answer-key updates.\n// Each table cell holds a unique '<dividend>\u00f7<divisor>=<quotient>, <remainder>' string,\n// so we locate each old value with body.search() and replace it in place,\n// which preserves the existing run/paragraph formatting (font, size, alignment).\nconst pairs = [\n  [\"71\u00f78=8, 7\", \"19\u00f76=3, 1\"],\n  [\"91\u00f77=13, 0\", \"34\u00f79=3, 7\"],\n  [\"56\u00f78=7, 0\", \"20\u00f78=2, 4\"],\n  [\"52\u00f76=8, 4\", \"74\u00f73=24, 2\"],\n  [\"86\u00f72=43, 0\", \"60\u00f73=20, 0\"],\n  [\"88\u00f74=22, 0\", \"15\u00f72=7, 1\"],\n  [\"11\u00f76=1, 5\", \"12\u00f79=1, 3\"],\n  [\"86\u00f73=28, 2\", \"29\u00f74=7, 1\"],\n  [\"21\u00f79=2, 3\", \"92\u00f76=15, 2\"],\n  [\"88\u00f79=9, 7\", \"69\u00f73=23, 0\"],\n  [\"71\u00f73=23, 2\", \"78\u00f73=26, 0\"],\n  [\"93\u00f74=23, 1\", \"58\u00f77=8, 2\"],\n  [\"33\u00f78=4, 1\", \"93\u00f75=18, 3\"],\n  [\"49\u00f77=7, 0\", \"27\u00f78=3, 3\"],\n  [\"90\u00f75=18, 0\", \"50\u00f72=25, 0\"],\n  [\"45\u00f72=22, 1\", \"23\u00f78=2, 7\"],\n  [\"85\u00f72=42, 1\", \"40\u00f76=6, 4\"],\n  [\"61\u00f78=7, 5\", \"47\u00f78=5, 7\"],\n  [\"70\u00f79=7, 7\", \"27\u00f77=3, 6\"],\n  [\"94\u00f75=18, 4\", \"53\u00f75=10, 3\"],\n  [\"92\u00f72=46, 0\", \"69\u00f76=11, 3\"],\n  [\"65\u00f72=32, 1\", \"15\u00f78=1, 7\"],\n  [\"78\u00f74=19, 2\", \"96\u00f74=24, 0\"],\n  [\"23\u00f77=3, 2\", \"85\u00f72=42, 1\"],\n  [\"98\u00f78=12, 2\", \"84\u00f78=10, 4\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly 1 match for '\" + oldText + \"' but found \" + results.items.length);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Apply the two-digit-division answer-key updates.\n# Each table cell holds a unique '<dividend>\u00f7<divisor>=<quotient>, <remainder>' string,\n# so we locate each old value with Find/Replace on the document range, which preserves\n# the existing run/paragraph formatting (font, size, alignment) of the matched text.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"71\u00f78=8, 7\", \"19\u00f76=3, 1\"),\n  @(\"91\u00f77=13, 0\", \"34\u00f79=3, 7\"),\n  @(\"56\u00f78=7, 0\", \"20\u00f78=2, 4\"),\n  @(\"52\u00f76=8, 4\", \"74\u00f73=24, 2\"),\n  @(\"86\u00f72=43, 0\", \"60\u00f73=20, 0\"),\n  @(\"88\u00f74=22, 0\", \"15\u00f72=7, 1\"),\n  @(\"11\u00f76=1, 5\", \"12\u00f79=1, 3\"),\n  @(\"86\u00f73=28, 2\", \"29\u00f74=7, 1\"),\n  @(\"21\u00f79=2, 3\", \"92\u00f76=15, 2\"),\n  @(\"88\u00f79=9, 7\", \"69\u00f73=23, 0\"),\n  @(\"71\u00f73=23, 2\", \"78\u00f73=26, 0\"),\n  @(\"93\u00f74=23, 1\", \"58\u00f77=8, 2\"),\n  @(\"33\u00f78=4, 1\", \"93\u00f75=18, 3\"),\n  @(\"49\u00f77=7, 0\", \"27\u00f78=3, 3\"),\n  @(\"90\u00f75=18, 0\", \"50\u00f72=25, 0\"),\n  @(\"45\u00f72=22, 1\", \"23\u00f78=2, 7\"),\n  @(\"85\u00f72=42, 1\", \"40\u00f76=6, 4\"),\n  @(\"61\u00f78=7, 5\", \"47\u00f78=5, 7\"),\n  @(\"70\u00f79=7, 7\", \"27\u00f77=3, 6\"),\n  @(\"94\u00f75=18, 4\", \"53\u00f75=10, 3\"),\n  @(\"92\u00f72=46, 0\", \"69\u00f76=11, 3\"),\n  @(\"65\u00f72=32, 1\", \"15\u00f78=1, 7\"),\n  @(\"78\u00f74=19, 2\", \"96\u00f74=24, 0\"),\n  @(\"23\u00f77=3, 2\", \"85\u00f72=42, 1\"),\n  @(\"98\u00f78=12, 2\", \"84\u00f78=10, 4\"),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $found = $find.Execute(\n    $oldText,  # FindText\n    $true,     # MatchCase\n    $false,    # MatchWholeWord\n    $false,    # MatchWildcards\n    $false,    # MatchSoundsLike\n    $false,    # MatchAllWordForms\n    $true,     # Forward\n    1,         # Wrap (wdFindContinue)\n    $false,    # Format\n    $newText,  # ReplaceWith\n    1          # Replace (wdReplaceOne)\n  )\n  if (-not $found) {\n    throw \"Could not find expected text: $oldText\"\n  }\n}\n"}
